# Update the RSS row (row 24) to reflect the new RSS colour.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B24: HEX colour code changes from FF6600 to FF8300
$ws.Range("B24").Value = "FF8300"

# C24: H value changes from 24 to 31 (D24 formula auto-recalculates to 131)
$ws.Range("C24").Value = 31

# Update the active selection to match the saved workbook state
$ws.Range("H26").Select()
